$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Thursday (row 16): In 08:00, Out 10:00
$ws.Range("C16").Value = 0.333333333333333
$ws.Range("D16").Value = 0.416666666666667

# Friday (row 17): In 08:00, Out 13:00 (split shift), with a note about the actual times worked
$ws.Range("C17").Value = 0.333333333333333
$ws.Range("D17").Value = 0.541666666666667
$ws.Range("L17").Value = "0800-1000, 1300 – 1500, 1800"

# Move the active selection to D18, matching the saved view state
$ws.Range("D18").Select()

# Re-apply the print area (recorded twice in the source workbook's history)
$ws.PageSetup.PrintArea = '$A$1:$K$27'
$ws.PageSetup.PrintArea = '$A$1:$K$27'
